# Update "想去人数" (interest/attendee counts) in column F across sheets,
# matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 8791
$ws1.Range("F24").Value = 2047
$ws1.Range("F27").Value = 1802
$ws1.Range("F33").Value = 114

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 10

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F14").Value = 8791
$ws4.Range("F26").Value = 2047
$ws4.Range("F29").Value = 1802
$ws4.Range("F35").Value = 114
$ws4.Range("F42").Value = 10
